$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Create the new "2022-Q1" sheet (holds the fund detail rows, same
#    layout as the other quarterly sheets) BEFORE moving it into place
#    -- the engine re-resolves worksheet handles by index on Move(),
#    so write all data first while $newQ still points at the right
#    physical sheet.
# ------------------------------------------------------------------
$newQ = $wb.Worksheets.Add()
$newQ.Name = "2022-Q1"

$newQ.Cells.Item(1,2).Value = "基金代码"
$newQ.Cells.Item(1,3).Value = "基金名称"
$newQ.Cells.Item(1,4).Value = "基金规模"
$newQ.Cells.Item(1,5).Value = "股票总仓位"
$newQ.Cells.Item(1,6).Value = "仓位占比"
$newQ.Cells.Item(1,7).Value = "持有市值(亿元)"
$newQ.Cells.Item(1,8).Value = "仓位排名"
$newQ.Range("B1:H1").Style = "Heading 1"

$newQ.Cells.Item(2,1).Value = 0

$newQ.Cells.Item(2,2).NumberFormat = "@"
$newQ.Cells.Item(2,2).Value = "006282"

$newQ.Cells.Item(2,3).Value = "上投摩根欧洲动力策略股票（QDII）"

$newQ.Cells.Item(2,4).NumberFormat = "@"
$newQ.Cells.Item(2,4).Value = "0.48"

$newQ.Cells.Item(2,5).NumberFormat = "@"
$newQ.Cells.Item(2,5).Value = "89.68"

$newQ.Cells.Item(2,6).NumberFormat = "@"
$newQ.Cells.Item(2,6).Value = "2.41"

$newQ.Cells.Item(2,7).NumberFormat = "@"
$newQ.Cells.Item(2,7).Value = "0.0116"

$newQ.Cells.Item(2,8).Value = 6

# Move "2022-Q1" so it sits right before "总计" (i.e. becomes the new
# second-to-last tab, pushing "总计" to the end).
$total = $wb.Worksheets.Item("总计")
$newQ.Move($total)

# ------------------------------------------------------------------
# 2) Update the "总计" (totals) sheet: insert a new row 2 for 2022-Q1,
#    pushing the previous rows down by one.
# ------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q1"
$total.Cells.Item(2,3).Value = 1
$total.Cells.Item(2,4).Value = 0.01

# Re-number the running index in column A (0,1,2,3,4,5) for every data
# row now that the table has grown by one row.
$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(4,1).Value = 2
$total.Cells.Item(5,1).Value = 3
$total.Cells.Item(6,1).Value = 4
$total.Cells.Item(7,1).Value = 5

# Row-insert carried stray formatting onto the new row (B2:D2 picked up
# the header's style, A2 lost its index-column style) -- fix both up so
# the layout matches the other rows in the table.
$total.Cells.Item(3,1).Copy()
$total.Cells.Item(2,1).PasteSpecial(-4122)
$total.Range("B2:D2").ClearFormats()
$excel.CutCopyMode = $false
